$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'63.439.96"
$ws.Range("E2").Value = "'  -2.11%  "
$ws.Range("D3").Value = "'3.158.03"
$ws.Range("E3").Value = "'  -0.26%  "
$ws.Range("E4").Value = "'  +0.05%  "
$ws.Range("D5").Value = "'591.81"
$ws.Range("E5").Value = "'  -2.25%  "
$ws.Range("D6").Value = "'139.11"
$ws.Range("E6").Value = "'  -4.46%  "
$ws.Range("E7").Value = "'  -0.06%  "
$ws.Range("D8").Value = "'3.152.63"
$ws.Range("E8").Value = "'  -0.10%  "
$ws.Range("E9").Value = "'  -1.57%  "
$ws.Range("E10").Value = "'  -2.80%  "
$ws.Range("D11").Value = "'5.28"
$ws.Range("E11").Value = "'  -2.37%  "
$ws.Range("E12").Value = "'  -2.96%  "
$ws.Range("D13").Value = "'0.0000247"
$ws.Range("E13").Value = "'  -3.90%  "
$ws.Range("D14").Value = "'34.29"
$ws.Range("E14").Value = "'  -4.05%  "
$ws.Range("D15").Value = "'3.677.33"
$ws.Range("E15").Value = "'  -0.11%  "
$ws.Range("D16").Value = "'0.121"
$ws.Range("E16").Value = "'  +1.17%  "
$ws.Range("D17").Value = "'3.150.01"
$ws.Range("E17").Value = "'  -1.01%  "
$ws.Range("D18").Value = "'63.357.21"
$ws.Range("E18").Value = "'  -2.23%  "
$ws.Range("D19").Value = "'6.69"
$ws.Range("E19").Value = "'  -3.45%  "
$ws.Range("D20").Value = "'477.28"
$ws.Range("E20").Value = "'  -1.41%  "
$ws.Range("D21").Value = "'14.15"
$ws.Range("E21").Value = "'  -4.09%  "
$ws.Range("D22").Value = "'0.702"
$ws.Range("E22").Value = "'  -2.11%  "
$ws.Range("E23").Value = "'  -0.33%  "
$ws.Range("D24").Value = "'84.52"
$ws.Range("E24").Value = "'  -3.80%  "
$ws.Range("D25").Value = "'13.05"
$ws.Range("E25").Value = "'  -4.07%  "
$ws.Range("E26").Value = "'  +0.00%  "
$ws.Range("E27").Value = "'  -2.23%  "
$ws.Range("B28").Value = "'NEARProtocol"
$ws.Range("C28").Value = "'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D28").Value = "'7.08"
$ws.Range("E28").Value = "'  -1.93%  "
$ws.Range("B29").Value = "'RenderToken"
$ws.Range("C29").Value = "'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D29").Value = "'8.04"
$ws.Range("E29").Value = "'  -5.59%  "
$ws.Range("E30").Value = "'  +1.33%  "
$ws.Range("E31").Value = "'  +0.05%  "
$ws.Range("D32").Value = "'26.93"
$ws.Range("E32").Value = "'  -1.11%  "
$ws.Range("E33").Value = "'  -5.36%  "
$ws.Range("D34").Value = "'2.55"
$ws.Range("E34").Value = "'  -6.16%  "
$ws.Range("E35").Value = "'  -3.02%  "
$ws.Range("E36").Value = "'  -4.51%  "
$ws.Range("D37").Value = "'52.57"
$ws.Range("E37").Value = "'  -0.92%  "
$ws.Range("D38").Value = "'0.0₃0705"
$ws.Range("E38").Value = "'  -8.53%  "
$ws.Range("D39").Value = "'0.0391"
$ws.Range("E39").Value = "'  -1.78%  "
$ws.Range("D40").Value = "'423.10"
$ws.Range("E40").Value = "'  -5.28%  "
$ws.Range("E41").Value = "'  -9.67%  "
$ws.Range("B42").Value = "'Cosmos"
$ws.Range("C42").Value = "'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D42").Value = "'8.28"
$ws.Range("E42").Value = "'  -0.44%  "
$ws.Range("B43").Value = "'Maker"
$ws.Range("C43").Value = "'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D43").Value = "'2.942.40"
$ws.Range("E43").Value = "'  +2.17%  "
$ws.Range("E44").Value = "'  -6.75%  "
$ws.Range("D45").Value = "'0.265"
$ws.Range("E45").Value = "'  -0.38%  "
$ws.Range("D46").Value = "'2.15"
$ws.Range("E46").Value = "'  -5.28%  "
$ws.Range("E47").Value = "'  +0.07%  "
$ws.Range("D48").Value = "'25.62"
$ws.Range("E48").Value = "'  -2.73%  "
$ws.Range("E49").Value = "'  -0.86%  "
$ws.Range("E50").Value = "'  -9.73%  "
$ws.Range("D51").Value = "'120.92"
$ws.Range("E51").Value = "'  -0.77%  "
